# Updating with proper kelp info and additional kelp data cleaning.
# The TRANSECT site "Norwegian" was renamed to "Norwegian Cove" in the
# "Data Entry" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Entry")

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value2 -eq "Norwegian") {
        $cell.Value2 = "Norwegian Cove"
    }
}
